$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 (Caso -239, NEWBERY JORGE) was removed from the source data.
# Rows 41-52 shift up to become rows 40-51; former row 52 is cleared.

# Force text storage (avoid auto numeric/date conversion) for columns A, B, E
$ws.Range("A40:B51").NumberFormat = "@"
$ws.Range("E40:E51").NumberFormat = "@"

# Row 40  (was row 41, Caso -241)
$ws.Range("A40").Value = "-241"
$ws.Range("B40").Value = "12/16/2024"
$ws.Range("C40").Value = "CUENCA /ALT/ 116"
$ws.Range("D40").Value = "106573 - FLORESTA"
$ws.Range("E40").Value = "801679775"
$ws.Range("F40").Value = "INCO"
$ws.Range("G40").Value = "Pendiente"
$ws.Range("H40").Value = "Podrida en base sin riego "
$ws.Range("I40").Value = "0"
$ws.Range("J40").Value = ""
$ws.Range("K40").Value = ""
$ws.Range("L40").Value = ""
$ws.Range("M40").Value = -58.474354
$ws.Range("N40").Value = -34.629997

# Row 41  (was row 42, Caso -247)
$ws.Range("A41").Value = "-247"
$ws.Range("B41").Value = "12/26/2024"
$ws.Range("C41").Value = "CONCORDIA /ALT/ 925"
$ws.Range("D41").Value = "106580 - SANTA RITA"
$ws.Range("E41").Value = "802055232"
$ws.Range("F41").Value = "INCO"
$ws.Range("G41").Value = "Pendiente"
$ws.Range("H41").Value = "Pasante con priroidad"
$ws.Range("I41").Value = "1"
$ws.Range("J41").Value = "Cambio"
$ws.Range("K41").Value = "Sin equipos"
$ws.Range("L41").Value = "Pasante"
$ws.Range("M41").Value = -58.479695
$ws.Range("N41").Value = -34.622867

# Row 42  (was row 43, Caso -251)
$ws.Range("A42").Value = "-251"
$ws.Range("B42").Value = "1/2/2025"
$ws.Range("C42").Value = "RINCON /ALT/ 645"
$ws.Range("D42").Value = "106558 - BALVANERA"
$ws.Range("E42").Value = "802269060"
$ws.Range("F42").Value = "INCO"
$ws.Range("G42").Value = "Pendiente"
$ws.Range("H42").Value = ""
$ws.Range("I42").Value = "1"
$ws.Range("J42").Value = "Cambio"
$ws.Range("K42").Value = "Sin equipos"
$ws.Range("L42").Value = "Pasante"
$ws.Range("M42").Value = -58.396131
$ws.Range("N42").Value = -34.616584

# Row 43  (was row 44, Caso -252)
$ws.Range("A43").Value = "-252"
$ws.Range("B43").Value = "1/2/2025"
$ws.Range("C43").Value = "LIBERTI TOMAS /ALT/ 1110"
$ws.Range("D43").Value = "106551 - LA BOCA"
$ws.Range("E43").Value = "802269071"
$ws.Range("F43").Value = "INCO"
$ws.Range("G43").Value = "Pendiente"
$ws.Range("H43").Value = ""
$ws.Range("I43").Value = "1"
$ws.Range("J43").Value = "Cambio"
$ws.Range("K43").Value = "Sin equipos"
$ws.Range("L43").Value = "Pasante"
$ws.Range("M43").Value = -58.369292
$ws.Range("N43").Value = -34.631678

# Row 44  (was row 45, Caso -255)
$ws.Range("A44").Value = "-255"
$ws.Range("B44").Value = "1/8/2025"
$ws.Range("C44").Value = "GURRUCHAGA /ALT/ 408"
$ws.Range("D44").Value = "106582 - VILLA CRESPO"
$ws.Range("E44").Value = "802393948"
$ws.Range("F44").Value = "INCO"
$ws.Range("G44").Value = "Pendiente"
$ws.Range("H44").Value = ""
$ws.Range("I44").Value = "1"
$ws.Range("J44").Value = "Cambio"
$ws.Range("K44").Value = "Nodo/Fuente Teco"
$ws.Range("L44").Value = "Pasante"
$ws.Range("M44").Value = -58.442667
$ws.Range("N44").Value = -34.597977

# Row 45  (was row 46, Caso -256)
$ws.Range("A45").Value = "-256"
$ws.Range("B45").Value = "1/8/2025"
$ws.Range("C45").Value = "NECOCHEA /ALT/ 1279"
$ws.Range("D45").Value = "106551 - LA BOCA"
$ws.Range("E45").Value = "802394092"
$ws.Range("F45").Value = "INCO"
$ws.Range("G45").Value = "Pendiente"
$ws.Range("H45").Value = "Picada"
$ws.Range("I45").Value = "1"
$ws.Range("J45").Value = "Cambio"
$ws.Range("K45").Value = "Sin equipos"
$ws.Range("L45").Value = "Pasante"
$ws.Range("M45").Value = -58.357221
$ws.Range("N45").Value = -34.635473

# Row 46  (was row 47, Caso -258)
$ws.Range("A46").Value = "-258"
$ws.Range("B46").Value = "1/14/2025"
$ws.Range("C46").Value = "CIUDAD DE LA PAZ /ALT/ 1465"
$ws.Range("D46").Value = "106581 - COLEGIALES"
$ws.Range("E46").Value = "802608477"
$ws.Range("F46").Value = "INCO"
$ws.Range("G46").Value = "Pendiente"
$ws.Range("H46").Value = "Picada"
$ws.Range("I46").Value = "1"
$ws.Range("J46").Value = "Cambio"
$ws.Range("K46").Value = "Sin equipos"
$ws.Range("L46").Value = "Pasante"
$ws.Range("M46").Value = -58.452317
$ws.Range("N46").Value = -34.567846

# Row 47  (was row 48, Caso -270)
$ws.Range("A47").Value = "-270"
$ws.Range("B47").Value = "1/27/2025"
$ws.Range("C47").Value = "SALTA SUR /ALT/ 917"
$ws.Range("D47").Value = "106552 - CONSTITUCION"
$ws.Range("E47").Value = "802925467"
$ws.Range("F47").Value = "INCO"
$ws.Range("G47").Value = "Pendiente"
$ws.Range("H47").Value = ""
$ws.Range("I47").Value = "1"
$ws.Range("J47").Value = "Cambio"
$ws.Range("K47").Value = "Sin equipos"
$ws.Range("L47").Value = "Pasante"
$ws.Range("M47").Value = -58.383027
$ws.Range("N47").Value = -34.618818

# Row 48  (was row 49, Caso -271)
$ws.Range("A48").Value = "-271"
$ws.Range("B48").Value = "1/27/2025"
$ws.Range("C48").Value = "HELGUERA /ALT/ 1405"
$ws.Range("D48").Value = "106580 - SANTA RITA"
$ws.Range("E48").Value = "802925468"
$ws.Range("F48").Value = "INCO"
$ws.Range("G48").Value = "Pendiente"
$ws.Range("H48").Value = "Base picada"
$ws.Range("I48").Value = "1"
$ws.Range("J48").Value = "Cambio"
$ws.Range("K48").Value = "Sin equipos"
$ws.Range("L48").Value = "Pasante"
$ws.Range("M48").Value = -58.480871
$ws.Range("N48").Value = -34.616598

# Row 49  (was row 50, Caso -275)
$ws.Range("A49").Value = "-275"
$ws.Range("B49").Value = "2/3/2025"
$ws.Range("C49").Value = "DEAN FUNES /ALT/ 481"
$ws.Range("D49").Value = "106558 - BALVANERA"
$ws.Range("E49").Value = "803039902"
$ws.Range("F49").Value = "INCO"
$ws.Range("G49").Value = "Pendiente"
$ws.Range("H49").Value = "Propia diámetro 114mm "
$ws.Range("I49").Value = "1"
$ws.Range("J49").Value = "Cambio"
$ws.Range("K49").Value = "Sin equipos"
$ws.Range("L49").Value = "Pasante"
$ws.Range("M49").Value = -58.407076
$ws.Range("N49").Value = -34.616016

# Row 50  (was row 51, Caso -282)
$ws.Range("A50").Value = "-282"
$ws.Range("B50").Value = "2/26/2025"
$ws.Range("C50").Value = "ALBARINO /ALT/ 1327"
$ws.Range("D50").Value = "106572 - MATADEROS"
$ws.Range("E50").Value = "803651213"
$ws.Range("F50").Value = "INCO"
$ws.Range("G50").Value = "Pendiente"
$ws.Range("H50").Value = "Podrida en la base"
$ws.Range("I50").Value = "0"
$ws.Range("J50").Value = "Cambio"
$ws.Range("K50").Value = "Sin equipos"
$ws.Range("L50").Value = "Pasante"
$ws.Range("M50").Value = -58.496341
$ws.Range("N50").Value = -34.650532

# Row 51  (was row 52, Caso -389)
$ws.Range("A51").Value = "-389"
$ws.Range("B51").Value = "5/4/2025"
$ws.Range("C51").Value = "AZARA /ALT/ 15"
$ws.Range("D51").Value = "106561 - BARRACAS"
$ws.Range("E51").Value = "805655333"
$ws.Range("F51").Value = "INCO"
$ws.Range("G51").Value = "Pendiente"
$ws.Range("H51").Value = ""
$ws.Range("I51").Value = "1"
$ws.Range("J51").Value = "Cambio"
$ws.Range("K51").Value = "Sin equipos"
$ws.Range("L51").Value = "Pasante"
$ws.Range("M51").Value = -58.372751
$ws.Range("N51").Value = -34.631917

# Remove the now-duplicate trailing row (former row 52)
$ws.Range("A52:N52").ClearContents()
